$d = $word.ActiveDocument

# Namespace / wrapper helpers for Range.InsertXML (WordprocessingML fragment package)
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rFonts = '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>'
$langEs = '<w:lang w:val="es-ES"/>'

# ---------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark from the title paragraph
#    (it will be re-added later at the end of the 2nd Results note)
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------
# 2) First "Results" notes paragraph: replace the English
#    placeholder text with the Spanish result text, and add the
#    es-ES language to the paragraph mark formatting too.
# ---------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range

$pPr5 = '<w:pPr><w:pStyle w:val="Notes"/><w:rPr>' + $rFonts + $langEs + '</w:rPr></w:pPr>'
$body5 = ""
$body5 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t>E</w:t></w:r>'
$body5 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t xml:space="preserve">l </w:t></w:r>'
$body5 += '<w:proofErr w:type="spellStart"/>'
$body5 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t>tester</w:t></w:r>'
$body5 += '<w:proofErr w:type="spellEnd"/>'
$body5 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t xml:space="preserve"> pudo detectar que existía un bug con el Nick, aunque no pudo determinar que </w:t></w:r>'
$body5 += '<w:proofErr w:type="spellStart"/>'
$body5 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t>AccountURL</w:t></w:r>'
$body5 += '<w:proofErr w:type="spellEnd"/>'
$body5 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t xml:space="preserve"> se mostraba dos veces.</w:t></w:r>'

$r5.InsertXML($pkgOpen + '<w:p>' + $pPr5 + $body5 + '</w:p>' + $pkgClose)

# ---------------------------------------------------------------
# 3) Second "Results" notes paragraph: replace the English
#    placeholder text with the Spanish result text, add the es-ES
#    language to the paragraph mark, and re-add the "_GoBack"
#    bookmark right at the end of the paragraph text.
# ---------------------------------------------------------------
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range

$pPr10 = '<w:pPr><w:pStyle w:val="Notes"/><w:rPr>' + $rFonts + $langEs + '</w:rPr></w:pPr>'
$body10 = ""
$body10 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t xml:space="preserve">El </w:t></w:r>'
$body10 += '<w:proofErr w:type="spellStart"/>'
$body10 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t>tester</w:t></w:r>'
$body10 += '<w:proofErr w:type="spellEnd"/>'
$body10 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t xml:space="preserve"> no pudo detectar este bug puesto que no accedió al listado de </w:t></w:r>'
$body10 += '<w:proofErr w:type="spellStart"/>'
$body10 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t>advertisement</w:t></w:r>'
$body10 += '<w:proofErr w:type="spellEnd"/>'
$body10 += '<w:r><w:rPr>' + $rFonts + $langEs + '</w:rPr><w:t xml:space="preserve"> en el display de ticket si no al listado de tickets de un usuario.</w:t></w:r>'
$body10 += '<w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/>'

$r10.InsertXML($pkgOpen + '<w:p>' + $pPr10 + $body10 + '</w:p>' + $pkgClose)

# ---------------------------------------------------------------
# 4) The two trailing empty paragraphs get the es-ES language on
#    their paragraph mark as well (no text runs, no pStyle).
# ---------------------------------------------------------------
$p11 = $d.Paragraphs(11)
$pPr11 = '<w:pPr><w:rPr>' + $rFonts + $langEs + '</w:rPr></w:pPr>'
$p11.Range.InsertXML($pkgOpen + '<w:p>' + $pPr11 + '</w:p>' + $pkgClose)

$p12 = $d.Paragraphs(12)
$pPr12 = '<w:pPr><w:rPr>' + $rFonts + $langEs + '</w:rPr></w:pPr>'
$p12.Range.InsertXML($pkgOpen + '<w:p>' + $pPr12 + '</w:p>' + $pkgClose)

Write-Host "Edit complete"
